$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update source data that changed (columns C, D, H, I, J for rows 6-14)
# ---------------------------------------------------------------------------

# Row 4 - reference Nifty close value
$ws.Range("J4").Value2 = 22404

# Row 6
$ws.Range("C6").Value2 = 11
$ws.Range("D6").Value2 = 8
$ws.Range("H6").Value2 = 13
$ws.Range("I6").Value2 = 7
$ws.Range("J6").Value2 = 22359

# Row 7
$ws.Range("C7").Value2 = 8
$ws.Range("D7").Value2 = 9
$ws.Range("H7").Value2 = 5
$ws.Range("I7").Value2 = 4
$ws.Range("J7").Value2 = 22400

# Row 8
$ws.Range("C8").Value2 = 5
$ws.Range("D8").Value2 = 8
$ws.Range("H8").Value2 = 4
$ws.Range("I8").Value2 = 2
$ws.Range("J8").Value2 = 22450

# Row 9
$ws.Range("C9").Value2 = 6
$ws.Range("D9").Value2 = 7
$ws.Range("H9").Value2 = 6
$ws.Range("I9").Value2 = 8
$ws.Range("J9").Value2 = 22469

# Row 10
$ws.Range("C10").Value2 = 9
$ws.Range("D10").Value2 = 7
$ws.Range("H10").Value2 = 6
$ws.Range("I10").Value2 = 7
$ws.Range("J10").Value2 = 22476

# Row 11
$ws.Range("C11").Value2 = 11
$ws.Range("D11").Value2 = 9
$ws.Range("H11").Value2 = 7
$ws.Range("I11").Value2 = 9
$ws.Range("J11").Value2 = 22475

# Row 12
$ws.Range("C12").Value2 = 37
$ws.Range("D12").Value2 = 34
$ws.Range("H12").Value2 = 29
$ws.Range("I12").Value2 = 21
$ws.Range("J12").Value2 = 22449

# Row 13 (new data - this row previously had no C/D/H/I/J values at all)
$ws.Range("C13").Value2 = 10
$ws.Range("D13").Value2 = 7
$ws.Range("H13").Value2 = 9
$ws.Range("I13").Value2 = 8
$ws.Range("J13").Value2 = 22474

# Row 14 (new data - this row previously had no C/D/H/I/J values at all)
$ws.Range("C14").Value2 = 28
$ws.Range("D14").Value2 = 36
$ws.Range("H14").Value2 = 34
$ws.Range("I14").Value2 = 27
$ws.Range("J14").Value2 = 22466

# ---------------------------------------------------------------------------
# 2. Resize Table3 (G5:K26 -> G5:L26) and add the new "DATA" column.
#    Set shared-string content in an order that reproduces the target
#    sharedStrings.xml append order: NIFTY, DATA, Time, ALL.
# ---------------------------------------------------------------------------

$lo = $ws.ListObjects.Item("Table3")
$lo.Resize($ws.Range("G5:L26"))

# L6:L11 - plain "NIFTY" labels (default/no explicit number format)
$ws.Range("L6").Value2 = "NIFTY"
$ws.Range("L7").Value2 = "NIFTY"
$ws.Range("L8").Value2 = "NIFTY"
$ws.Range("L9").Value2 = "NIFTY"
$ws.Range("L10").Value2 = "NIFTY"
$ws.Range("L11").Value2 = "NIFTY"

# L5 - new table column header "DATA"
$ws.Range("L5").Value2 = "DATA"

# V5 - new "Time" header, matching the style of the other bold table headers
$ws.Range("M5").Copy()
$ws.Range("V5").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("V5").Value2 = "Time"

# L12 - "ALL" label
$ws.Range("L12").Value2 = "ALL"

# L13/L14 carry the numeric (0.00) style inherited from column K, plus text
$ws.Range("K4").Copy()
$ws.Range("L13").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("L13").Value2 = "NIFTY"

$ws.Range("K4").Copy()
$ws.Range("L14").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("L14").Value2 = "ALL"

# ---------------------------------------------------------------------------
# 3. New "Time" data column (V7:V14), formatted as time-of-day (h:mm)
# ---------------------------------------------------------------------------

$ws.Range("V7:V14").NumberFormat = "h:mm"
$ws.Range("V7").Value2 = 0.42777777777777781
$ws.Range("V8").Value2 = 0.45208333333333334
$ws.Range("V9").Value2 = 0.46736111111111112
$ws.Range("V10").Value2 = 0.57986111111111105
$ws.Range("V11").Value2 = 0.59652777777777777
$ws.Range("V12").Value2 = 0.62986111111111109
$ws.Range("V13").Value2 = 0.6381944444444444
$ws.Range("V14").Value2 = 0.65

# ---------------------------------------------------------------------------
# 4. Sheet view: update selection (closest achievable approximation of the
#    saved view state - active cell X20, scrolled near row 3).
# ---------------------------------------------------------------------------

$win = $excel.Windows.Item(1)
$win.ScrollRow = 3
$win.ScrollColumn = 1
$ws.Range("X20").Select()
